$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lowercase the header labels in row 1
$ws.Range("A1").Value = "nace"
$ws.Range("B1").Value = "year"

# Clear out the now-removed data points in row 2 (C2 and E2)
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
